# Automatische test-sync: 2025-06-30 20:15:50
#
# Adds Testmail #17 (newsletter unsubscribe request) to the "Logs" sheet,
# swaps the "Bestelling / Levering" / "Openingstijden / Locatie" rows on
# the "Dashboard" sheet back to their correct order, appends a new
# "Afmelding / Nieuwsbrief" tally row, and extends the conditional
# formatting / chart series ranges to cover the grown data.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 17 with the new test mail entry.
# ---------------------------------------------------------------------
$logs.Cells.Item(17, 1).Value = "Kunt u mij uitschrijven voor de nieuwsbrief?"
$logs.Cells.Item(17, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(17, 3).Value = "Testmail #17: Kunt u mij uitschrijven voor de nieuwsbrief?"
$logs.Cells.Item(17, 4).Value = "Afmelding / Nieuwsbrief"

$logs.Cells.Item(17, 5).Value = @"
Beste klant,
Dank voor uw bericht. Om u uit te schrijven voor onze nieuwsbrief, hebben wij uw e-mailadres nodig. Kunt u ons alstublieft het e-mailadres sturen waar u voor uitgeschreven wilt worden?
Met vriendelijke groet,
[Bedrijfsnaam] E-mailassistent
"@

$logs.Cells.Item(17, 6).Value = "2025-06-30 20:14:54"
$logs.Cells.Item(17, 7).Value = "Ja"
$logs.Cells.Item(17, 8).Value = "Nee"
$logs.Cells.Item(17, 9).Value = "Ja"
$logs.Cells.Item(17, 10).Value = "Nee"

# The multi-line text just written into column E makes the engine
# auto-expand the row height (mirrors real Excel's implicit re-wrap on
# newline-containing values). None of the other rows carry an explicit
# height, so re-run AutoFit to drop the row back to an un-pinned
# (non-custom) standard height, matching the source file.
$logs.Rows.Item(17).AutoFit()

# ---------------------------------------------------------------------
# 2. Logs sheet: extend the conditional formatting ranges from row 16
#    to row 17 for every formatted column.
# ---------------------------------------------------------------------
foreach ($col in "D", "G", "H", "I", "J") {
    $oldRange = $logs.Range("$($col)2:$($col)16")
    $newRange = $logs.Range("$($col)2:$($col)17")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 3. Dashboard sheet: rows 5 and 6 had their categories swapped;
#    restore "Openingstijden / Locatie" to row 5 and
#    "Bestelling / Levering" to row 6 (counts stay at 1 each).
# ---------------------------------------------------------------------
$dash.Cells.Item(5, 1).Value = "Openingstijden / Locatie"
$dash.Cells.Item(6, 1).Value = "Bestelling / Levering"

# ---------------------------------------------------------------------
# 4. Dashboard sheet: append the new "Afmelding / Nieuwsbrief" tally row.
# ---------------------------------------------------------------------
$dash.Cells.Item(10, 1).Value = "Afmelding / Nieuwsbrief"
$dash.Cells.Item(10, 2).Value = 1

# ---------------------------------------------------------------------
# 5. Chart: grow the category/value series references from row 9 to
#    row 10 to include the new Dashboard row.
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$10,Dashboard!`$B`$2:`$B`$10,1)"
